$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[0.21342326297962266, 13.00082627543468]"
$ws.Range("N2").Value = 0.04312273696604918
$ws.Range("O2").Value = 0.04312273696604918
$ws.Range("Q2").Value = "[-2.817684702318773, -0.050315798255692545]"
$ws.Range("R2").Value = 0.04255030281703265
$ws.Range("S2").Value = 0.04255030281703265
$ws.Range("U2").Value = "[5.018990335976898, 12.673312238931985]"
$ws.Range("V2").Value = 0.00002869836718466523
$ws.Range("W2").Value = 0.00002869836718466523
$ws.Range("Y2").Value = 0.2081281281281342
$ws.Range("Z2").Value = 11.65517517517545
$ws.Range("M3").Value = "[-1.0899303542357472, 13.858948977358665]"
$ws.Range("N3").Value = 0.09222960119982915
$ws.Range("O3").Value = 0.09222960119982915
$ws.Range("Q3").Value = "[-3.547263777026312, -0.1509473947670772]"
$ws.Range("R3").Value = 0.03350684771698154
$ws.Range("S3").Value = 0.03350684771698154
$ws.Range("U3").Value = "[5.461435903498645, 13.878063875870277]"
$ws.Range("V3").Value = 0.0000313946745831295
$ws.Range("W3").Value = 0.0000313946745831295
$ws.Range("Y3").Value = 0.6243843843844008
$ws.Range("Z3").Value = 14.67303303303338
$ws.Range("M4").Value = "[-1.4766008015267928, 14.490480815988356]"
$ws.Range("N4").Value = 0.1076479270738941
$ws.Range("O4").Value = 0.1076479270738941
$ws.Range("Q4").Value = "[-5.4466851611787, 0.2138421425866932]"
$ws.Range("R4").Value = 0.06915149521547592
$ws.Range("S4").Value = 0.06915149521547592
$ws.Range("U4").Value = "[5.163313059760018, 13.668525241460195]"
$ws.Range("V4").Value = 0.00005425782262769019
$ws.Range("W4").Value = 0.00005425782262769019
$ws.Range("Y4").Value = -0.8845445445445677
$ws.Range("Z4").Value = 22.52986986987041
$ws.Range("M5").Value = "[-1.5567959394952737, 14.656355068060062]"
$ws.Range("N5").Value = 0.1106535720877817
$ws.Range("O5").Value = 0.1106535720877817
$ws.Range("Q5").Value = "[-4.9686850777496225, -0.4654211338651546]"
$ws.Range("R5").Value = 0.0191310717017763
$ws.Range("S5").Value = 0.0191310717017763
$ws.Range("U5").Value = "[5.201939489069154, 13.646531952292921]"
$ws.Range("V5").Value = 0.0000482971888462469
$ws.Range("W5").Value = 0.0000482971888462469
$ws.Range("Y5").Value = 1.925185185185232
$ws.Range("Z5").Value = 20.55265265265314
$ws.Range("M6").Value = "[-1.1054528029386965, 13.9346716645165]"
$ws.Range("N6").Value = 0.0926659751991703
$ws.Range("O6").Value = 0.0926659751991703
$ws.Range("Q6").Value = "[1.0755001877154236, 4.647921863869584]"
$ws.Range("R6").Value = 0.002336210289328067
$ws.Range("S6").Value = 0.002336210289328067
$ws.Range("U6").Value = "[5.300293617419893, 13.737570054748051]"
$ws.Range("V6").Value = 0.00004118535529640255
$ws.Range("W6").Value = 0.00004118535529640255
$ws.Range("Y6").Value = 6.764164164164326
$ws.Range("Z6").Value = 21.54126126126178
$ws.Range("M7").Value = "[-1.3119892586116482, 14.588723765999704]"
$ws.Range("N7").Value = 0.09954936626885647
$ws.Range("O7").Value = 0.09954936626885647
$ws.Range("Q7").Value = "[-0.4780000834290776, 4.943527178621776]"
$ws.Range("R7").Value = 0.1040811733444098
$ws.Range("S7").Value = 0.1040811733444098
$ws.Range("U7").Value = "[5.343658130559504, 13.530686813158177]"
$ws.Range("V7").Value = 0.00002985854158832346
$ws.Range("W7").Value = 0.00002985854158832346
$ws.Range("Y7").Value = 5.285555555555653
$ws.Range("Z7").Value = 26.6759259259264
$ws.Range("M8").Value = "[-0.6167033346194675, 14.443063505172956]"
$ws.Range("N8").Value = 0.0710120130011549
$ws.Range("O8").Value = 0.0710120130011549
$ws.Range("Q8").Value = "[0.2704474156243464, 3.3648690083494266]"
$ws.Range("R8").Value = 0.02233955324453185
$ws.Range("S8").Value = 0.02233955324453185
$ws.Range("U8").Value = "[5.452075339439258, 13.632287455338236]"
$ws.Range("V8").Value = 0.00002489065378319388
$ws.Range("W8").Value = 0.00002489065378319388
$ws.Range("Y8").Value = 11.51407407407428
$ws.Range("Z8").Value = 23.72296296296338
$ws.Range("M9").Value = "[-1.1973998942430004, 14.601561705649049]"
$ws.Range("N9").Value = 0.09437637953349753
$ws.Range("O9").Value = 0.09437637953349753
$ws.Range("Q9").Value = "[0.11950002085726918, 3.9938164865455823]"
$ws.Range("R9").Value = 0.03795292781593385
$ws.Range("S9").Value = 0.03795292781593385
$ws.Range("U9").Value = "[5.349312466936334, 13.550228595600002]"
$ws.Range("V9").Value = 0.0000300224620077838
$ws.Range("W9").Value = 0.0000300224620077838
$ws.Range("Y9").Value = 9.032592592592756
$ws.Range("Z9").Value = 24.31851851851895
$ws.Range("B10").Value = 1
$ws.Range("M10").Value = "[0.07054734953545427, 13.643466606433797]"
$ws.Range("N10").Value = 0.04776278397956446
$ws.Range("O10").Value = 0.04776278397956446
$ws.Range("Q10").Value = "[0.42139481039142357, 3.4906585039886595]"
$ws.Range("R10").Value = 0.0136487085209076
$ws.Range("S10").Value = 0.0136487085209076
$ws.Range("U10").Value = "[4.719046483838294, 12.262689132753344]"
$ws.Range("V10").Value = 0.00004262907435670371
$ws.Range("W10").Value = 0.00004262907435670371
$ws.Range("Y10").Value = 11.01777777777797
$ws.Range("Z10").Value = 23.12740740740781
$ws.Range("M11").Value = "[0.032125852029606605, 13.243417470181045]"
$ws.Range("N11").Value = 0.04894226675991264
$ws.Range("O11").Value = 0.04894226675991264
$ws.Range("Q11").Value = "[0.5094474573388856, 3.050395269251349]"
$ws.Range("R11").Value = 0.007084844960868697
$ws.Range("S11").Value = 0.007084844960868697
$ws.Range("U11").Value = "[4.223141229267657, 11.549261142003663]"
$ws.Range("V11").Value = 0.00008068398959193424
$ws.Range("W11").Value = 0.00008068398959193424
$ws.Range("Y11").Value = 12.75481481481505
$ws.Range("Z11").Value = 22.7800000000004
$ws.Range("M12").Value = "[-0.7282631646868936, 13.496990411718148]"
$ws.Range("N12").Value = 0.07731181292783873
$ws.Range("O12").Value = 0.07731181292783873
$ws.Range("Q12").Value = "[-0.42139481039142446, 2.8994478744842738]"
$ws.Range("R12").Value = 0.1398383151045526
$ws.Range("S12").Value = 0.1398383151045526
$ws.Range("U12").Value = "[4.4637375903534675, 12.27884319244911]"
$ws.Range("V12").Value = 0.00008636732020028681
$ws.Range("W12").Value = 0.00008636732020028681
$ws.Range("Y12").Value = 13.3503703703706
$ws.Range("Z12").Value = 26.45259259259306
